# Update metadata sheet to reflect suppressed OIDs:
#  - "Experimental" row now carries a value ("false")
#  - "Date" row is refreshed to the new generation timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B7").Value = "'false"
$ws.Range("B8").Value = "2023-10-09T22:41:16+02:00"
